$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) replacing the old Strike# values, rows 2-18
$kValues = @(6, 1, 4, 4, 3, 7, 4, 0, 3, 4, 1, 2, 0, 5, 4, 1, 4)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
